$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update 최종점수 (K) and MACRO_SCORE (N) columns for rows 2-7
$ws.Range("K2").Value = 61.9
$ws.Range("N2").Value = 49.16024380385575

$ws.Range("K3").Value = 53.7
$ws.Range("N3").Value = 49.16024380385575

$ws.Range("K4").Value = 49.9
$ws.Range("N4").Value = 49.16024380385575

$ws.Range("K5").Value = 47.1
$ws.Range("N5").Value = 49.16024380385575

$ws.Range("K6").Value = 39.1
$ws.Range("N6").Value = 49.16024380385575

$ws.Range("K7").Value = 39.1
$ws.Range("N7").Value = 49.16024380385575

$wb.Save()
